$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark price cells that look like plain numbers as Text so Excel keeps the
# original string (with trailing zeros / leading zeros) instead of coercing
# them to a floating point number.
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D22","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '27.588.01'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').Value = '1.853.55'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('D4').Value = '1.032'
$ws.Range('E4').Value = '  +2.62%  '
$ws.Range('D5').Value = '322.36'
$ws.Range('E5').Value = '  +3.46%  '
$ws.Range('D6').Value = '1.028'
$ws.Range('E6').Value = '  +2.33%  '
$ws.Range('D7').Value = '0.4394'
$ws.Range('E7').Value = '  +2.51%  '
$ws.Range('D8').Value = '0.3787'
$ws.Range('E8').Value = '  +2.74%  '
$ws.Range('D9').Value = '0.07414'
$ws.Range('E9').Value = '  +2.69%  '
$ws.Range('D10').Value = '0.8792'
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('D11').Value = '21.67'
$ws.Range('E11').Value = '  +2.55%  '
$ws.Range('D12').Value = '1.866.08'
$ws.Range('E12').Value = '  -7.98%  '
$ws.Range('D13').Value = '5.528'
$ws.Range('E13').Value = '  +2.76%  '
$ws.Range('D14').Value = '6.708'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').Value = '0.07208'
$ws.Range('E15').Value = '  +4.61%  '
$ws.Range('D16').Value = '83.26'
$ws.Range('E16').Value = '  +3.15%  '
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '0.000009070'
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').Value = '1.029'
$ws.Range('E19').Value = '  +2.40%  '
$ws.Range('D20').Value = '15.46'
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('D21').Value = '27.618.65'
$ws.Range('E21').Value = '  +2.37%  '
$ws.Range('D22').Value = '5.284'
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('E23').Value = '  +4.21%  '
$ws.Range('D24').Value = '158.01'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('D25').Value = '1.913'
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('D26').Value = '18.77'
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('D27').Value = '1.984'
$ws.Range('E27').Value = '  +4.81%  '
$ws.Range('D28').Value = '5.292'
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').Value = '117.45'
$ws.Range('E29').Value = '  +2.23%  '
$ws.Range('D30').Value = '0.09065'
$ws.Range('E30').Value = '  +1.87%  '
$ws.Range('D31').Value = '1.205'
$ws.Range('E31').Value = '  +4.42%  '
$ws.Range('D32').Value = '0.7629'
$ws.Range('E32').Value = '  +2.93%  '
$ws.Range('D33').Value = '4.546'
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('D34').Value = '2.885'
$ws.Range('E34').Value = '  +2.94%  '
$ws.Range('D35').Value = '1.030'
$ws.Range('E35').Value = '  +1.95%  '
$ws.Range('D36').Value = '1.152'
$ws.Range('E36').Value = '  +3.14%  '
$ws.Range('D37').Value = '0.01980'
$ws.Range('E37').Value = '  +3.19%  '
$ws.Range('D38').Value = '0.05317'
$ws.Range('E38').Value = '  +2.21%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.5173'
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.823'
$ws.Range('E40').Value = '  +2.98%  '
$ws.Range('E41').Value = '  +2.43%  '
$ws.Range('D42').Value = '6.773'
$ws.Range('E42').Value = '  +5.67%  '
$ws.Range('D43').Value = '8.586'
$ws.Range('E43').Value = '  +4.25%  '
$ws.Range('D44').Value = '108.96'
$ws.Range('E44').Value = '  +2.10%  '
$ws.Range('D45').Value = '10.55'
$ws.Range('E45').Value = '  +2.31%  '
$ws.Range('D46').Value = '1.716'
$ws.Range('E46').Value = '  +4.16%  '
$ws.Range('D47').Value = '0.4660'
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('D48').Value = '0.06403'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('E49').Value = '  +3.27%  '
$ws.Range('D50').Value = '39.46'
$ws.Range('E50').Value = '  +4.71%  '
$ws.Range('D51').Value = '64.12'
$ws.Range('E51').Value = '  +0.89%  '
